$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "subject" header to "subject d"
$ws.Range("C1").Value = "subject d"

# Move the active selection to C1 (matches the cell that was just edited)
$ws.Range("C1").Select()
